$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the career name: "Informáticas" -> "Informáticos"
$ws.Range("A9").Value = "Ingeniería de Sistemas Informáticos"

# Update the active selection to A9, matching the edited cell
$ws.Range("A9").Select()
